$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the "POPULAR" programming-language table ---
# Row 9 used to duplicate row 7 ("Python" / 0.388); it now holds the
# "C#" figure that used to live two rows further down.
$ws.Range("B9").Value = "C#"
$ws.Range("C9").Value = 0.34399999999999997

# Row 10 now holds what used to be row 11's PHP figure (STT jumps 5 -> 7).
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "PHP"
$ws.Range("C10").Value = 0.307

# Row 11 now holds what used to be row 12's Ruby figure.
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Ruby"
$ws.Range("C11").Value = 0.10100000000000001

# Row 12 is now the trailing (now unused) row of the table - blank it out
# but keep its formatting.
$ws.Range("A12:C12").ClearContents()

# The spacer row 14 becomes completely empty (value + formatting), so it
# collapses out of the sheet's used range entirely.
$ws.Range("B14").Clear()

# Row 15 gains a matching, still-empty C15 cell (same percent formatting
# as the rest of the spacer column).
$ws.Range("C15").NumberFormat = "0.00%"

# Move the active selection from E16 to B16.
$ws.Range("B16").Select()

# --- Keep the pie chart's source range in sync with the shrunken table ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Sheet1!`$C`$4,Sheet1!`$B`$5:`$B`$11,Sheet1!`$C`$5:`$C`$11,1)"

# Best-effort: nudge the plot area into the same manually-laid-out
# position/size (expressed as fractions of the chart area, edge-anchored)
# that Excel recorded after the data range shrank.
$plotArea = $chart.PlotArea
$plotArea.Left = 0.22130005200208008
$plotArea.Top = 0.11155988594391172
$plotArea.Width = 0.56572006502307337
$plotArea.Height = 0.78746258489288035
